$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'259.10"
$ws.Range("E2").Value = "'0.68%"
$ws.Range("D3").Value = "'27.04"
$ws.Range("E3").Value = "'-0.31%"
$ws.Range("D4").Value = "'4.704"
$ws.Range("E4").Value = "'0.72%"
$ws.Range("D5").Value = "'0.06021"
$ws.Range("E5").Value = "'2.29%"
$ws.Range("D6").Value = "'6.675"
$ws.Range("E6").Value = "'0.44%"
$ws.Range("D7").Value = "'0.8600"
$ws.Range("E7").Value = "'0.33%"
$ws.Range("D8").Value = "'0.9273"
$ws.Range("E8").Value = "'-4.11%"
$ws.Range("D9").Value = "'0.1397"
$ws.Range("E9").Value = "'-0.83%"
$ws.Range("D10").Value = "'0.04977"
$ws.Range("E10").Value = "'24.77%"
$ws.Range("D11").Value = "'0.07030"
$ws.Range("D12").Value = "'0.03092"
$ws.Range("E12").Value = "'-2.80%"
$ws.Range("D13").Value = "'0.09125"
$ws.Range("E13").Value = "'-0.51%"
$ws.Range("E14").Value = "'-0.32%"
$ws.Range("D15").Value = "'0.0006046"
$ws.Range("E15").Value = "'-0.40%"
$ws.Range("D16").Value = "'0.005988"
$ws.Range("E16").Value = "'-3.70%"
$ws.Range("D17").Value = "'3.465"
$ws.Range("E17").Value = "'-1.47%"
$ws.Range("D18").Value = "'3.169"
$ws.Range("E18").Value = "'-1.04%"
$ws.Range("E19").Value = "'-1.83%"
$ws.Range("E20").Value = "'0.43%"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("E21").Value = "'0.44%"
$ws.Range("D22").Value = "'4.132"
$ws.Range("E22").Value = "'6.78%"
$ws.Range("D23").Value = "'0.04246"
$ws.Range("E23").Value = "'0.68%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'-0.41%"
$ws.Range("D25").Value = "'0.004036"
$ws.Range("E26").Value = "'-0.06%"
$ws.Range("D27").Value = "'0.0001523"
$ws.Range("E27").Value = "'-21.36%"
$ws.Range("D40").Value = "'0.03844"
$ws.Range("E40").Value = "'0.37%"
$ws.Range("E41").Value = "'1.01%"
$ws.Range("D42").Value = "'0.003988"
$ws.Range("E42").Value = "'1.73%"
$ws.Range("E43").Value = "'31.24%"
$ws.Range("E44").Value = "'-9.53%"
$ws.Range("D45").Value = "'0.00005101"
$ws.Range("E45").Value = "'-6.58%"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.1321"
$ws.Range("E47").Value = "'0.48%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.05454"
$ws.Range("E48").Value = "'-9.10%"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E50").Value = "'-0.01%"
